$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.699.16'
$ws.Range("E2").Value = '  -2.83%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.744.02'
$ws.Range("E3").Value = '  -4.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.69'
$ws.Range("E5").Value = '  -8.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5055'
$ws.Range("E7").Value = '  -5.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.74'
$ws.Range("E8").Value = '  -6.86%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2651'
$ws.Range("E9").Value = '  -11.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06124'
$ws.Range("E10").Value = '  -10.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.746.33'
$ws.Range("E11").Value = '  -5.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06959'
$ws.Range("E12").Value = '  -4.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.23'
$ws.Range("E13").Value = '  -14.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.489'
$ws.Range("E14").Value = '  -9.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5963'
$ws.Range("E15").Value = '  -18.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '76.62'
$ws.Range("E16").Value = '  -13.82%  '
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.699.92'
$ws.Range("E19").Value = '  -2.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006800'
$ws.Range("E20").Value = '  -13.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.60'
$ws.Range("E21").Value = '  -16.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.969.23'
$ws.Range("E22").Value = '  -5.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.035'
$ws.Range("E23").Value = '  -11.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.142'
$ws.Range("E24").Value = '  -11.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.096'
$ws.Range("E25").Value = '  -14.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.40'
$ws.Range("E26").Value = '  -3.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.516'
$ws.Range("E27").Value = '  -9.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.823'
$ws.Range("E28").Value = '  -17.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.96'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '103.22'
$ws.Range("E30").Value = '  -6.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.758'
$ws.Range("E31").Value = '  -11.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08104'
$ws.Range("E32").Value = '  -7.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.459'
$ws.Range("E33").Value = '  -13.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04491'
$ws.Range("E34").Value = '  -6.32%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.653'
$ws.Range("E36").Value = '  -9.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9823'
$ws.Range("E37").Value = '  -13.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6086'
$ws.Range("E38").Value = '  -16.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.651'
$ws.Range("E39").Value = '  -14.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01550'
$ws.Range("E40").Value = '  -9.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.921'
$ws.Range("E41").Value = '  -16.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.53'
$ws.Range("E43").Value = '  -3.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3799'
$ws.Range("E44").Value = '  -19.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.119'
$ws.Range("E45").Value = '  -12.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7273'
$ws.Range("E46").Value = '  -19.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05332'
$ws.Range("E47").Value = '  -7.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1109'
$ws.Range("E48").Value = '  -9.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.08'
$ws.Range("E49").Value = '  -13.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.893'
$ws.Range("E50").Value = '  -19.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.39'
$ws.Range("E51").Value = '  -12.82%  '
